$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice: rebuild the Osm-Lifr LR-pair result table (rows 2-19)
# Clear old rows 2-13 first so stale cells beyond the new table are also handled.
$ws.Range("A2:T13").ClearContents()

$colA = @("M1", "M1", "M1", "M1", "M1", "M1", "M2", "M2", "M2", "M2", "M2", "M2", "Neutro", "Neutro", "Neutro", "Neutro", "Neutro", "Neutro")
$colB = @("Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm", "Osm")
$colC = @("Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr", "Lifr")
$colD = @("ECs", "FAPs", "M1", "M2", "Neutro", "sCs", "ECs", "FAPs", "M1", "M2", "Neutro", "sCs", "ECs", "FAPs", "M1", "M2", "Neutro", "sCs")
$colE = @(3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 3.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
$colF = @(1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333, 0.3333333333333333)
$colG = @(26.18311833333333, 26.18311833333333, 26.18311833333333, 26.18311833333333, 26.18311833333333, 26.18311833333333, 31.31438633333333, 31.31438633333333, 31.31438633333333, 31.31438633333333, 31.31438633333333, 31.31438633333333, 0.1591176666666667, 0.1591176666666667, 0.1591176666666667, 0.1591176666666667, 0.1591176666666667, 0.1591176666666667)
$colH = @(78.54935499999999, 78.54935499999999, 78.54935499999999, 78.54935499999999, 78.54935499999999, 78.54935499999999, 93.943159, 93.943159, 93.943159, 93.943159, 93.943159, 93.943159, 0.477353, 0.477353, 0.477353, 0.477353, 0.477353, 0.477353)
$colI = @(0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853)
$colJ = @(0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.4541216129859197, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.5431186404276995, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853, 0.002759746586380853)
$colK = @(2.0, 3.0, 3.0, 3.0, 3.0, 2.0, 2.0, 3.0, 3.0, 3.0, 3.0, 2.0, 2.0, 3.0, 3.0, 3.0, 3.0, 2.0)
$colL = @(1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0, 1.0)
$colM = @(31.7631495, 17.81777966666667, 13.32746566666667, 24.04189266666667, 23.53771866666667, 16.857219, 31.7631495, 17.81777966666667, 13.32746566666667, 24.04189266666667, 23.53771866666667, 16.857219, 31.7631495, 17.81777966666667, 13.32746566666667, 24.04189266666667, 23.53771866666667, 16.857219)
$colN = @(63.52629899999999, 53.453339, 39.982397, 72.12567800000001, 70.613156, 33.714438, 63.52629899999999, 53.453339, 39.982397, 72.12567800000001, 70.613156, 33.714438, 63.52629899999999, 53.453339, 39.982397, 72.12567800000001, 70.613156, 33.714438)
$colO = @(0.2494255238736205, 0.139917139754138, 0.1046561867492399, 0.1887930437533183, 0.184833931824778, 0.1323741740449054, 0.2494255238736205, 0.139917139754138, 0.1046561867492399, 0.1887930437533183, 0.184833931824778, 0.1323741740449054, 0.2494255238736205, 0.139917139754138, 0.1046561867492399, 0.1887930437533183, 0.184833931824778, 0.1323741740449054)
$colP = @(0.1905320411699034, 0.1603205907999899, 0.1199177007191215, 0.2163238354260682, 0.211787385034485, 0.1011184468504321, 0.1905320411699034, 0.1603205907999899, 0.1199177007191215, 0.2163238354260682, 0.211787385034485, 0.1011184468504321, 0.1905320411699034, 0.1603205907999899, 0.1199177007191215, 0.2163238354260682, 0.211787385034485, 0.1011184468504321)
$colQ = @(831.6583019978573, 466.5250334495938, 348.9546106337705, 629.4917206486322, 616.2908731460421, 441.374559847915, 994.6435346064234, 557.9528360842112, 417.3414087302359, 752.8571151485336, 737.0692157333115, 527.873468271607, 5.0540782344245, 2.835123525740778, 2.120635239460111, 3.825489863370445, 3.745266872896445, 2.682281353769)
$colR = @(4989.949811987144, 4198.725301046345, 3140.591495703934, 5665.42548583769, 5546.61785831438, 2648.24735908749, 5967.86120763854, 5021.575524757901, 3756.072678572123, 6775.714036336803, 6633.622941599804, 3167.240809629642, 30.324469406547, 25.516111731667, 19.085717155141, 34.42940877033401, 33.707401856068, 16.093688122614)
$colS = @(0.1132695212213466, 0.06353939718952548, 0.04752663633552044, 0.08573500154977821, 0.08393708325479769, 0.06011397343495131, 0.1354676514142075, 0.07599160671579984, 0.05684072585959458, 0.1025370212455094, 0.1003867537575795, 0.07189488143500869, 0.0006883512380664801, 0.0003861358488126549, 0.0002888245541248518, 0.0005210209580306711, 0.0005100948124007823, 0.0003653191749454126)
$colT = @(0.08652471786157621, 0.07280504528894702, 0.05445721967613022, 0.09823732907098676, 0.09617722890193035, 0.04592007218634922, 0.1034815031581124, 0.08707310130785606, 0.06512953857778501, 0.1174895073887116, 0.1150256766196672, 0.05491931337556728, 0.0005258201502148171, 0.0004424442031868336, 0.0003309424652062362, 0.0005969989663699052, 0.0005844795128874472, 0.0002790612885156137)

for ($i = 0; $i -lt 18; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $colA[$i]
    $ws.Cells.Item($r, 2).Value2 = $colB[$i]
    $ws.Cells.Item($r, 3).Value2 = $colC[$i]
    $ws.Cells.Item($r, 4).Value2 = $colD[$i]
    $ws.Cells.Item($r, 5).Value2 = $colE[$i]
    $ws.Cells.Item($r, 6).Value2 = $colF[$i]
    $ws.Cells.Item($r, 7).Value2 = $colG[$i]
    $ws.Cells.Item($r, 8).Value2 = $colH[$i]
    $ws.Cells.Item($r, 9).Value2 = $colI[$i]
    $ws.Cells.Item($r, 10).Value2 = $colJ[$i]
    $ws.Cells.Item($r, 11).Value2 = $colK[$i]
    $ws.Cells.Item($r, 12).Value2 = $colL[$i]
    $ws.Cells.Item($r, 13).Value2 = $colM[$i]
    $ws.Cells.Item($r, 14).Value2 = $colN[$i]
    $ws.Cells.Item($r, 15).Value2 = $colO[$i]
    $ws.Cells.Item($r, 16).Value2 = $colP[$i]
    $ws.Cells.Item($r, 17).Value2 = $colQ[$i]
    $ws.Cells.Item($r, 18).Value2 = $colR[$i]
    $ws.Cells.Item($r, 19).Value2 = $colS[$i]
    $ws.Cells.Item($r, 20).Value2 = $colT[$i]
}
